# "remove duplicate records for today all"
#
# The sheet holds a rolling log of notification checks (one row per date).
# This edit prepends 10 newer daily records (2018-05-08 .. 2018-03-19) above
# the two rows that were already present (2018-03-16, 2018-03-15), which end
# up shifted down to rows 12-13. The final sheet has 12 data rows (rows 2-13)
# plus the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row order for rows 2..13 (date, col C, col D, col E); column F is
# always "Y" and column A is always 0 (styled like the existing index column).
$dates  = @("2018-05-08","2018-05-07","2018-05-04","2018-05-03","2018-05-02","2018-03-23","2018-03-22","2018-03-21","2018-03-20","2018-03-19","2018-03-16","2018-03-15")
$cVals  = @(0.6910412695946914,0.4990321303156645,0.2083083340506761,0.4245870500400741,0.459373184580592,0.6897610196508391,0.008445144199579573,0.04849964445249828,0.1240769640212653,0.4128998839088854,0.5094264856692802,0.4862442484091188)
$dVals  = @(0.6615130147022792,0.4824120289985133,0.2052835283182502,0.4631555243935099,0.4856046565429817,0.7198472011673879,0.1745463354600005,0.1456137080423014,0.1258930767857917,0.3857020776683011,0.4788491921005814,0.4619932283919192)
$eVals  = @(0.6812175883782291,0.5233755845902032,0.1528411347936059,0.3318268654883984,0.3753580884130054,0.762582669176173,0.3724702861574418,0.3839333625597994,0.3388044588952401,0.5513699538613661,0.6570775583089501,0.5978473548984882)

# Column A (rows 4-13 are brand new cells) needs the same direct formatting
# (border/bold/center) that A2/A3 already use. Copy that formatting down
# instead of re-creating it, so no new cell style gets introduced.
$ws.Range("A2").Copy()
$ws.Range("A4:A13").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Column B holds plain text dates like "2018-05-08". If we just assign the
# string, Excel auto-recognizes it as a date and stores a serial number
# instead of text, so force the whole destination range to Text format first.
$ws.Range("B2:B13").NumberFormat = "@"

for ($i = 0; $i -lt 12; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = 0
    $ws.Cells.Item($row, 2).Value = $dates[$i]
    $ws.Cells.Item($row, 3).Value = $cVals[$i]
    $ws.Cells.Item($row, 4).Value = $dVals[$i]
    $ws.Cells.Item($row, 5).Value = $eVals[$i]
    $ws.Cells.Item($row, 6).Value = "Y"
}

# Clear the temporary Text formatting back to the default (unstyled) look
# that column B originally had, now that the literal text is stored.
$ws.Range("B2:B13").Style = "Normal"
